$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9998639178507085
$ws.Range("E2").Value = 0.9998639178507085

# Row 3
$ws.Range("D3").Value = 0.9999999958228338
$ws.Range("E3").Value = 0.9999999958228338

# Row 4
$ws.Range("D4").Value = 0.3962371857946776
$ws.Range("E4").Value = 0.3962371857946776

# Row 5
$ws.Range("D5").Value = 0.9999999998775795
$ws.Range("E5").Value = 0.9999999998775795

# Row 6
$ws.Range("D6").Value = 0.9784905462750674
$ws.Range("E6").Value = 0.9784905462750674

# Row 7
$ws.Range("D7").Value = 1.0
$ws.Range("E7").Value = 0.0

# Row 8
$ws.Range("D8").Value = 0.999941112858042
$ws.Range("E8").Value = 0.00005888714195800837

# Row 9
$ws.Range("D9").Value = 0.6805352703188723
$ws.Range("E9").Value = 0.3194647296811277

# Row 10
$ws.Range("D10").Value = 0.9981100231905022
$ws.Range("E10").Value = 0.001889976809497762

# Row 11
$ws.Range("D11").Value = 0.9993452489908626
$ws.Range("E11").Value = 0.0006547510091373665
$ws.Range("F11").Value = 5.575076103210449

# Row 12
$ws.Range("D12").Value = 0.9999879413270308
$ws.Range("E12").Value = 0.9999879413270308

# Row 13
$ws.Range("D13").Value = 0.9999997477762933
$ws.Range("E13").Value = 0.9999997477762933

# Row 14
$ws.Range("D14").Value = 0.0416629642831565
$ws.Range("E14").Value = 0.0416629642831565

# Row 15
$ws.Range("D15").Value = 0.9941997522073657
$ws.Range("E15").Value = 0.9941997522073657

# Row 16
$ws.Range("D16").Value = 0.999490337616115
$ws.Range("E16").Value = 0.999490337616115

# Row 17
$ws.Range("D17").Value = 1.0
$ws.Range("E17").Value = 0.0

# Row 18
$ws.Range("D18").Value = 0.9999999971968725
$ws.Range("E18").Value = 0.000000002803127507888803

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.00000000101707360045382
$ws.Range("E19").Value = 0.9999999989829264

# Row 20
$ws.Range("D20").Value = 0.9999759883274187
$ws.Range("E20").Value = 0.00002401167258125891

# Row 21
$ws.Range("D21").Value = 0.9999983992336746
$ws.Range("E21").Value = 0.00000160076632538253
$ws.Range("F21").Value = 5.999920845031738
$ws.Range("G21").Value = 0.5

Write-Output "done"
